$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: Trade #2 closed -> capital / P&L / trade counters bump up
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1300.07               # Current Capital
$summary.Range("B4").Value = 0.07000000000000001   # Total P&L $
$summary.Range("B5").Value = 0.7                   # Total P&L %
$summary.Range("B6").Value = 2                     # Total Trades
$summary.Range("B7").Value = 2                     # Winning Trades

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row picks up the new closed trade
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.07                 # Capital
$status.Range("D4").Value = 2                      # Trades
$status.Range("E4").Value = 0.07000000000000001    # P&L $
$status.Range("F4").Value = 0.07000000000000001    # P&L %

# ---------------------------------------------------------------------------
# Helper: write a trade-log row (#2) identically onto a trade-history sheet.
# Date/Time columns must stay plain text (not get coerced into date/time
# serials by Excel's usual type inference), so force a text number format
# before assigning, then clear the leftover formatting so the cell ends up
# with the default style, same as every other text cell on the sheet.
# ---------------------------------------------------------------------------
function Write-TradeRow2($ws) {
    $ws.Cells.Item(3, 1).Value = 2                 # Trade #

    $ws.Cells.Item(3, 2).NumberFormat = "@"
    $ws.Cells.Item(3, 2).Value = "2026-02-17"       # Date
    $ws.Cells.Item(3, 2).ClearFormats()

    $ws.Cells.Item(3, 3).NumberFormat = "@"
    $ws.Cells.Item(3, 3).Value = "19:42:35"         # Time
    $ws.Cells.Item(3, 3).ClearFormats()

    $ws.Cells.Item(3, 4).Value = "MarketMaking"     # Strategy
    $ws.Cells.Item(3, 5).Value = "DOWN"             # Side
    $ws.Cells.Item(3, 6).Value = 0.34               # Entry Price
    $ws.Cells.Item(3, 7).Value = 0.4                # Exit Price
    $ws.Cells.Item(3, 8).Value = "CLOSED"           # Status
    $ws.Cells.Item(3, 9).Value = 17.6471            # P&L %
    $ws.Cells.Item(3, 10).Value = 0.06              # P&L $
    $ws.Cells.Item(3, 11).Value = 100.07            # Capital After
    $ws.Cells.Item(3, 12).Value = 0                 # Entry Slippage (bps)
    $ws.Cells.Item(3, 13).Value = 0                 # Exit Slippage (bps)
    $ws.Cells.Item(3, 14).Value = 0.6               # Confidence
    $ws.Cells.Item(3, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(3, 16).Value = "early_exit"      # Exit Reason
    $ws.Cells.Item(3, 17).Value = 0.13              # Duration (min)
}

Write-TradeRow2 $wb.Worksheets.Item("All Trades")
Write-TradeRow2 $wb.Worksheets.Item("MarketMaking")
